# Alex Carey dagger sheet: insert ownTeam/oppTeam columns, re-order match
# rows (newest first) and append the Sharjah match that was previously
# missing.
#
# The sheet stores every value as plain TEXT, even ones that look numeric
# (e.g. "14", "107.69"). For a cell outside the sheet's original used
# range, Excel auto-types a numeric-looking string as a real number unless
# the cell is pre-formatted as Text ("@") - so that format is applied right
# before writing any such value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header)
$ws.Cells.Item(1, 1).Value = "venue"
$ws.Cells.Item(1, 2).Value = "date"
$ws.Cells.Item(1, 3).Value = "result"
$ws.Cells.Item(1, 4).Value = "ownTeam"
$ws.Cells.Item(1, 5).Value = "oppTeam"
$ws.Cells.Item(1, 6).Value = "batsman"
$ws.Cells.Item(1, 7).Value = "totalRuns"
$ws.Cells.Item(1, 8).Value = "totalBalls"
$ws.Cells.Item(1, 9).Value = "total4s"
$ws.Cells.Item(1, 10).Value = "total6s"
$ws.Cells.Item(1, 11).Value = "sr"

# Row 2 (Dubai (DSC) vs Rajasthan Royals)
$ws.Cells.Item(2, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(2, 2).Value = " October 14 2020"
$ws.Cells.Item(2, 3).Value = "Capitals won by 13 runs"
$ws.Cells.Item(2, 4).Value = "Delhi Capitals"
$ws.Cells.Item(2, 5).Value = "Rajasthan Royals"
$ws.Cells.Item(2, 6).Value = "Alex Carey †"
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "14"
$ws.Cells.Item(2, 8).NumberFormat = "@"
$ws.Cells.Item(2, 8).Value = "13"
$ws.Cells.Item(2, 9).NumberFormat = "@"
$ws.Cells.Item(2, 9).Value = "0"
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "1"
$ws.Cells.Item(2, 11).NumberFormat = "@"
$ws.Cells.Item(2, 11).Value = "107.69"

# Row 3 (Sharjah vs Chennai Super Kings)
$ws.Cells.Item(3, 1).Value = " Sharjah"
$ws.Cells.Item(3, 2).Value = " October 17 2020"
$ws.Cells.Item(3, 3).Value = "Capitals won by 5 wickets (with 1 ball remaining)"
$ws.Cells.Item(3, 4).Value = "Delhi Capitals"
$ws.Cells.Item(3, 5).Value = "Chennai Super Kings"
$ws.Cells.Item(3, 6).Value = "Alex Carey †"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "4"
$ws.Cells.Item(3, 8).NumberFormat = "@"
$ws.Cells.Item(3, 8).Value = "7"
$ws.Cells.Item(3, 9).NumberFormat = "@"
$ws.Cells.Item(3, 9).Value = "0"
$ws.Cells.Item(3, 10).NumberFormat = "@"
$ws.Cells.Item(3, 10).Value = "0"
$ws.Cells.Item(3, 11).NumberFormat = "@"
$ws.Cells.Item(3, 11).Value = "57.14"

# Row 4 (Abu Dhabi vs Mumbai Indians)
$ws.Cells.Item(4, 1).Value = " Abu Dhabi"
$ws.Cells.Item(4, 2).Value = " October 11 2020"
$ws.Cells.Item(4, 3).Value = "Mumbai won by 5 wickets (with 2 balls remaining)"
$ws.Cells.Item(4, 4).Value = "Delhi Capitals"
$ws.Cells.Item(4, 5).Value = "Mumbai Indians"
$ws.Cells.Item(4, 6).Value = "Alex Carey †"
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "14"
$ws.Cells.Item(4, 8).NumberFormat = "@"
$ws.Cells.Item(4, 8).Value = "9"
$ws.Cells.Item(4, 9).NumberFormat = "@"
$ws.Cells.Item(4, 9).Value = "0"
$ws.Cells.Item(4, 10).NumberFormat = "@"
$ws.Cells.Item(4, 10).Value = "0"
$ws.Cells.Item(4, 11).NumberFormat = "@"
$ws.Cells.Item(4, 11).Value = "155.55"
